$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("events")
$ws.Activate()

# Remove the "comment" column (column C) from the events sheet.
# This shifts the old column D (event_date_or_datetime) into column C.
$ws.Columns("C").Delete()

# Update selection to reflect the new active cell after the edit.
$ws.Range("F3").Select()
